# Collapse the 4 "gbParcelsShipmentService..." rows (117-120) into a single
# row describing the generic delivery-type selector, then drop the stray
# "Submit" button-selector value that trails the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 117 in place to the new, combined entry...
$ws.Range("A117").Value = "gbParcelsShipmentServiceDeliveryType"
$ws.Range("B117").Value = "input[id*='generalSelectedPickupDeliveryOptionCode'] + label[for='%s']"

# ...then delete the 3 rows that used to hold the per-option variants
# (118:120), which shifts everything below up by three rows.
$ws.Rows("118:120").Delete()

# The former last row (Submit) loses its paired CSS-selector value; after
# the shift above it now lives at row 160. Reset style before clearing so
# the now-empty cell is dropped entirely rather than left as a styled
# placeholder.
$ws.Range("B160").Style = "Normal"
$ws.Range("B160").ClearContents()
